# Apply the BOM rework: split the 3.9K and 10K resistor groups into
# "STD" (populated) and "DNP" (do-not-populate) rows, which pushes the
# IC rows (U1..U4) down by two rows, then restore the print layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two blank rows before the U1 row (old row 18) -------------
# This shifts old rows 18-21 (U1, U2, U3, U4) down to rows 20-23 and
# extends the shared "=Bn*5" formula fill range automatically.
$ws.Rows("18:19").Insert()

# --- 2. Fill the new rows' "Reference(s)" column first (bottom row up),
#        matching the order the BOM was actually re-split in: the 10K
#        DNP remainder, then the 10K STD remainder, then the 3.9K DNP
#        remainder -- this is also what the new rows need in full below.
$ws.Cells.Item(19, 4).Value = "R8, R13, R22"
$ws.Cells.Item(18, 4).Value = "R7, R9, R10, R11, R12, R14, R23, R24, R28"
$ws.Cells.Item(16, 4).Value = "R16, R17"

# --- 3. Row 19 (new): the 10K "DNP" remainder (R8, R13, R22) -------------
$ws.Cells.Item(19, 2).Value = 3
$ws.Cells.Item(19, 3).Formula = "=B19*5"
$ws.Cells.Item(19, 5).Value = "RC0805JR-0710KL"
$ws.Cells.Item(19, 6).Value = "DNP"
$ws.Cells.Item(19, 7).Value = "10K"
$ws.Cells.Item(19, 8).Value = "Resistor_SMD:R_0805_2012Metric"

# --- 4. Row 18 (new): the 10K "STD" group, minus R8/R13/R22 --------------
$ws.Cells.Item(18, 2).Value = 12
$ws.Cells.Item(18, 3).Formula = "=B18*5"
$ws.Cells.Item(18, 5).Value = "RC0805JR-0710KL"
$ws.Cells.Item(18, 6).Value = "STD"
$ws.Cells.Item(18, 7).Value = "10K"
$ws.Cells.Item(18, 8).Value = "Resistor_SMD:R_0805_2012Metric"

# --- 5. Row 16: was the 3.9K "STD" group (R3,R4,R16,R17) -> becomes the
#        DNP-only remainder of that group (R16, R17) ---------------------
$ws.Cells.Item(16, 2).Value = 2
$ws.Cells.Item(16, 5).Value = "RC0805FR-073K9L"
$ws.Cells.Item(16, 6).Value = "DNP"
$ws.Cells.Item(16, 7).Value = "3.9K"
$ws.Cells.Item(16, 8).Value = "Resistor_SMD:R_0805_2012Metric"

# --- 6. Row 17: becomes the old 100R group (previously row 16) -----------
$ws.Cells.Item(17, 2).Value = 10
$ws.Cells.Item(17, 4).Value = "R5, R6, R15, R18, R19, R20, R21, R25, R26, R27"
$ws.Cells.Item(17, 5).Value = "RC0805FR-07100RL"
$ws.Cells.Item(17, 6).Value = "STD"
$ws.Cells.Item(17, 7).Value = "100R"
$ws.Cells.Item(17, 8).Value = "Resistor_SMD:R_0805_2012Metric"

# --- 7. Re-number the "Item" column (col A) for the rows below the split -
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(23, 1).Value = 22

# --- 8. Recalc so the "=Bn*5" totals (col C) pick up the new quantities --
$ws.Calculate()

# --- 9. Sheet view: selection moves to D28, dimension grows to A1:H23 ---
[void]$ws.Range("D28").Select()

# --- 10. Print setup: explicit print area, narrower L/R margins, landscape
$ws.PageSetup.PrintArea = "A1:H23"
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.25)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.25)
$ws.PageSetup.Orientation = 2
